$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 40
$ws1.Range("F4").Value = 257
$ws1.Range("F6").Value = 195
$ws1.Range("F7").Value = 259
$ws1.Range("F8").Value = 50
$ws1.Range("F11").Value = 45
$ws1.Range("F12").Value = 29
$ws1.Range("F13").Value = 95
$ws1.Range("F14").Value = 1978
$ws1.Range("F15").Value = 50
$ws1.Range("F16").Value = 13
$ws1.Range("F17").Value = 496
$ws1.Range("F18").Value = 460
$ws1.Range("F22").Value = 40
$ws1.Range("F23").Value = 1476
$ws1.Range("F24").Value = 3393
$ws1.Range("F28").Value = 1098
$ws1.Range("F29").Value = 86
$ws1.Range("F30").Value = 1794
$ws1.Range("F33").Value = 58
$ws1.Range("F34").Value = 277
$ws1.Range("F37").Value = 643
$ws1.Range("F39").Value = 42

# 演出 (Performances) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 9

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 40
$ws4.Range("F4").Value = 257
$ws4.Range("F6").Value = 195
$ws4.Range("F7").Value = 259
$ws4.Range("F8").Value = 50
$ws4.Range("F11").Value = 45
$ws4.Range("F12").Value = 29
$ws4.Range("F13").Value = 95
$ws4.Range("F14").Value = 1978
$ws4.Range("F15").Value = 50
$ws4.Range("F16").Value = 9
$ws4.Range("F17").Value = 13
$ws4.Range("F18").Value = 496
$ws4.Range("F19").Value = 460
$ws4.Range("F23").Value = 40
$ws4.Range("F24").Value = 1476
$ws4.Range("F25").Value = 3393
$ws4.Range("F29").Value = 1098
$ws4.Range("F30").Value = 87
$ws4.Range("F31").Value = 1794
$ws4.Range("F34").Value = 58
$ws4.Range("F35").Value = 277
$ws4.Range("F38").Value = 643
$ws4.Range("F40").Value = 42
